# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-18 on Sheet1 to reflect the recalculated K statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 3
    4  = 3
    5  = 0
    6  = 2
    7  = 0
    8  = 2
    9  = 3
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
